$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 43: "2035_TM152_FBP_Plus_01" / FinalBlueprint / Plus ---
$ws.Rows(43).Insert()
$ws.Range("A42:H42").Copy()
$ws.Range("A43:H43").PasteSpecial(-4122)
$ws.Range("A43").Value = "RTP2021"
$ws.Range("B43").Value = 2035
$ws.Range("C43").Value = "2035_TM152_FBP_Plus_01"
$ws.Range("D43").Value = "FinalBlueprint"
$ws.Range("E43").Value = "Plus"
$ws.Range("F43").Value = '"Blueprint Plus Crossing (s23)\v1.7.1- FINAL DRAFT BLUEPRINT"'
$ws.Range("G43").Value = "run98"
$ws.Range("H43").Value = "current"

# --- Append new row 62: "2050_TM152_FBP_PlusCrossing_01" / FinalBlueprint / Plus ---
$ws.Range("A44:H44").Copy()
$ws.Range("A62:H62").PasteSpecial(-4122)
$ws.Range("A62").Value = "RTP2021"
$ws.Range("B62").Value = 2050
$ws.Range("C62").Value = "2050_TM152_FBP_PlusCrossing_01"
$ws.Range("D62").Value = "FinalBlueprint"
$ws.Range("E62").Value = "Plus"
$ws.Range("F62").Value = '"Blueprint Plus Crossing (s23)\v1.7.1- FINAL DRAFT BLUEPRINT"'
$ws.Range("G62").Value = "run98"
$ws.Range("H62").Value = "current"
